$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Ben Smit" typo in the per-person label above the Team-Total
# column (N5) so it matches the already-existing "Ben Smith" string used
# elsewhere in the sheet. Once nothing references the "Ben Smit" shared
# string any more it drops out of sharedStrings.xml entirely.
$ws.Range("N5").Value = "Ben Smith"

# Log a new activity entry in the 4th (J:L) block, row 10 -- "Added some
# time" per the commit message. Copy the date formatting from the cell
# above (J9) first so the new date cell keeps the same date number format
# / border style as the rest of the column instead of the plain default.
$ws.Range("J9").Copy()
$ws.Range("J10").PasteSpecial(-4122)

$ws.Range("J10").Value = 42608
$ws.Range("K10").Value = "Looked over Docs"
$ws.Range("L10").Value = 0.25

# Reflect the author's final cursor position.
$ws.Range("K17").Select() | Out-Null
